$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank first row, shifting the table (originally A2:F16) up to A1:F15
$ws.Rows.Item(1).Delete()

# Update the active selection to match the post-edit state
$ws.Range("B2").Select()
